$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 3
    "G2" = 241.2872163333334
    "H2" = 723.8616490000001
    "I2" = 0.8139476876902633
    "J2" = 0.8139476876902633
    "K2" = 3
    "M2" = 229.5846506666667
    "N2" = 688.753952
    "O2" = 0.5033187360873315
    "P2" = 0.5033187360873315
    "Q2" = 55395.84127222077
    "R2" = 498562.5714499869
    "S2" = 0.4096751214094693
    "T2" = 0.4096751214094693
    "E3" = 3
    "G3" = 241.2872163333334
    "H3" = 723.8616490000001
    "I3" = 0.8139476876902633
    "J3" = 0.8139476876902633
    "K3" = 3
    "M3" = 135.7283196666666
    "N3" = 407.1849589999999
    "O3" = 0.2975573763642838
    "P3" = 0.2975573763642838
    "Q3" = 32749.50842997082
    "R3" = 294745.5758697374
    "S3" = 0.2421961384468902
    "T3" = 0.2421961384468902
    "E4" = 3
    "G4" = 241.2872163333334
    "H4" = 723.8616490000001
    "I4" = 0.8139476876902633
    "J4" = 0.8139476876902633
    "K4" = 3
    "M4" = 90.23148833333335
    "N4" = 270.694465
    "O4" = 0.1978146123067711
    "P4" = 0.1978146123067711
    "Q4" = 21771.70464556365
    "R4" = 195945.3418100728
    "S4" = 0.1610107462784422
    "T4" = 0.1610107462784422
    "E5" = 3
    "G5" = 241.2872163333334
    "H5" = 723.8616490000001
    "I5" = 0.8139476876902633
    "J5" = 0.8139476876902633
    "K5" = 3
    "M5" = 0.5972149999999999
    "N5" = 1.791645
    "O5" = 0.001309275241613694
    "P5" = 0.001309275241613694
    "Q5" = 144.1003449025117
    "R5" = 1296.903104122605
    "S5" = 0.001065681555461577
    "T5" = 0.001065681555461577
    "E6" = 3
    "G6" = 42.36118633333333
    "H6" = 127.083559
    "I6" = 0.1428993636206566
    "J6" = 0.1428993636206566
    "K6" = 3
    "M6" = 229.5846506666667
    "N6" = 688.753952
    "O6" = 0.5033187360873315
    "P6" = 0.5033187360873315
    "Q6" = 9725.478166163906
    "R6" = 87529.30349547516
    "S6" = 0.07192392708523289
    "T6" = 0.07192392708523289
    "E7" = 3
    "G7" = 42.36118633333333
    "H7" = 127.083559
    "I7" = 0.1428993636206566
    "J7" = 0.1428993636206566
    "K7" = 3
    "M7" = 135.7283196666666
    "N7" = 407.1849589999999
    "O7" = 0.2975573763642838
    "P7" = 0.2975573763642838
    "Q7" = 5749.612640109896
    "R7" = 51746.51376098907
    "S7" = 0.04252075972308837
    "T7" = 0.04252075972308837
    "E8" = 3
    "G8" = 42.36118633333333
    "H8" = 127.083559
    "I8" = 0.1428993636206566
    "J8" = 0.1428993636206566
    "K8" = 3
    "M8" = 90.23148833333335
    "N8" = 270.694465
    "O8" = 0.1978146123067711
    "P8" = 0.1978146123067711
    "Q8" = 3822.312890422326
    "R8" = 34400.81601380094
    "S8" = 0.0282675822135045
    "T8" = 0.0282675822135045
    "E9" = 3
    "G9" = 42.36118633333333
    "H9" = 127.083559
    "I9" = 0.1428993636206566
    "J9" = 0.1428993636206566
    "K9" = 3
    "M9" = 0.5972149999999999
    "N9" = 1.791645
    "O9" = 0.001309275241613694
    "P9" = 0.001309275241613694
    "Q9" = 25.29873589606166
    "R9" = 227.688623064555
    "S9" = 0.0001870945988308784
    "T9" = 0.0001870945988308784
    "E10" = 3
    "G10" = 1.281292333333333
    "H10" = 3.843877
    "I10" = 0.004322255226862814
    "J10" = 0.004322255226862814
    "K10" = 3
    "M10" = 229.5846506666667
    "N10" = 688.753952
    "O10" = 0.5033187360873315
    "P10" = 0.5033187360873315
    "Q10" = 294.1650527502115
    "R10" = 2647.485474751904
    "S10" = 0.002175472037831454
    "T10" = 0.002175472037831454
    "E11" = 3
    "G11" = 1.281292333333333
    "H11" = 3.843877
    "I11" = 0.004322255226862814
    "J11" = 0.004322255226862814
    "K11" = 3
    "M11" = 135.7283196666666
    "N11" = 407.1849589999999
    "O11" = 0.2975573763642838
    "P11" = 0.2975573763642838
    "Q11" = 173.9076554051158
    "R11" = 1565.168898646043
    "S11" = 0.001286118925282111
    "T11" = 0.001286118925282111
    "E12" = 3
    "G12" = 1.281292333333333
    "H12" = 3.843877
    "I12" = 0.004322255226862814
    "J12" = 0.004322255226862814
    "K12" = 3
    "M12" = 90.23148833333335
    "N12" = 270.694465
    "O12" = 0.1978146123067711
    "P12" = 0.1978146123067711
    "Q12" = 115.6129142267561
    "R12" = 1040.516228040805
    "S12" = 0.0008550052419927823
    "T12" = 0.0008550052419927823
    "E13" = 3
    "G13" = 1.281292333333333
    "H13" = 3.843877
    "I13" = 0.004322255226862814
    "J13" = 0.004322255226862814
    "K13" = 3
    "M13" = 0.5972149999999999
    "N13" = 1.791645
    "O13" = 0.001309275241613694
    "P13" = 0.001309275241613694
    "Q13" = 0.7652070008516666
    "R13" = 6.886863007665
    "S13" = 0.000005659021756466864
    "T13" = 0.000005659021756466864
    "E14" = 3
    "G14" = 1.328211
    "H14" = 3.984633
    "I14" = 0.004480528594276053
    "J14" = 0.004480528594276054
    "K14" = 3
    "M14" = 229.5846506666667
    "N14" = 688.753952
    "O14" = 0.5033187360873315
    "P14" = 0.5033187360873315
    "Q14" = 304.9368584466239
    "R14" = 2744.431726019616
    "S14" = 0.002255133989074171
    "T14" = 0.002255133989074171
    "E15" = 3
    "G15" = 1.328211
    "H15" = 3.984633
    "I15" = 0.004480528594276053
    "J15" = 0.004480528594276054
    "K15" = 3
    "M15" = 135.7283196666666
    "N15" = 407.1849589999999
    "O15" = 0.2975573763642838
    "P15" = 0.2975573763642838
    "Q15" = 180.2758471927829
    "R15" = 1622.482624735047
    "S15" = 0.001333214333237935
    "T15" = 0.001333214333237935
    "E16" = 3
    "G16" = 1.328211
    "H16" = 3.984633
    "I16" = 0.004480528594276053
    "J16" = 0.004480528594276054
    "K16" = 3
    "M16" = 90.23148833333335
    "N16" = 270.694465
    "O16" = 0.1978146123067711
    "P16" = 0.1978146123067711
    "Q16" = 119.846455350705
    "R16" = 1078.618098156345
    "S16" = 0.0008863140268061193
    "T16" = 0.0008863140268061195
    "E17" = 3
    "G17" = 1.328211
    "H17" = 3.984633
    "I17" = 0.004480528594276053
    "J17" = 0.004480528594276054
    "K17" = 3
    "M17" = 0.5972149999999999
    "N17" = 1.791645
    "O17" = 0.001309275241613694
    "P17" = 0.001309275241613694
    "Q17" = 0.7932275323649998
    "R17" = 7.139047791284999
    "S17" = 0.000005866245157827846
    "T17" = 0.000005866245157827847
    "E18" = 3
    "G18" = 10.18278666666667
    "H18" = 30.54836
    "I18" = 0.03435016486794112
    "J18" = 0.03435016486794112
    "K18" = 3
    "M18" = 229.5846506666667
    "N18" = 688.753952
    "O18" = 0.5033187360873315
    "P18" = 0.5033187360873315
    "Q18" = 2337.811519679858
    "R18" = 21040.30367711872
    "S18" = 0.01728908156572358
    "T18" = 0.01728908156572358
    "E19" = 3
    "G19" = 10.18278666666667
    "H19" = 30.54836
    "I19" = 0.03435016486794112
    "J19" = 0.03435016486794112
    "K19" = 3
    "M19" = 135.7283196666666
    "N19" = 407.1849589999999
    "O19" = 0.2975573763642838
    "P19" = 0.2975573763642838
    "Q19" = 1382.092523790804
    "R19" = 12438.83271411724
    "S19" = 0.01022114493578515
    "T19" = 0.01022114493578515
    "E20" = 3
    "G20" = 10.18278666666667
    "H20" = 30.54836
    "I20" = 0.03435016486794112
    "J20" = 0.03435016486794112
    "K20" = 3
    "M20" = 90.23148833333335
    "N20" = 270.694465
    "O20" = 0.1978146123067711
    "P20" = 0.1978146123067711
    "Q20" = 918.8079963141557
    "R20" = 8269.2719668274
    "S20" = 0.006794964546025441
    "T20" = 0.006794964546025441
    "E21" = 3
    "G21" = 10.18278666666667
    "H21" = 30.54836
    "I21" = 0.03435016486794112
    "J21" = 0.03435016486794112
    "K21" = 3
    "M21" = 0.5972149999999999
    "N21" = 1.791645
    "O21" = 0.001309275241613694
    "P21" = 0.001309275241613694
    "Q21" = 6.081312939133332
    "R21" = 54.73181645219999
    "S21" = 0.00004497382040694385
    "T21" = 0.00004497382040694385
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output ("Updated " + $updates.Count + " cells")